# Weekly fruit/vegetable price update: insert a new data row above the
# existing last row (row 20), shifting that row down to row 21, and
# populate the new row 20 with the updated weekly record
# (Fecha 2021-11-09 / Volumen 100) while the previous record
# (Fecha 2021-10-22 / Volumen 60) survives unchanged as row 21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 20, pushing the old
# row 20 down to row 21 (this also extends the used dimension to R21).
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly record.
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 44509
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112026
$ws.Range("G20").Value = "Haba"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = 8500
$ws.Range("N20").Value = "$/saco 25 kilos"
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 340
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"
